$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 268
$ws.Range("A268").Value = "21TRD09437"
$ws.Range("B268").Value = "Bunner"
$ws.Range("C268").Value = "DUS - AMENDED to Distracted Driving"
$ws.Range("D268").Value = "'4510.11"
$ws.Range("E268").Value = "M1"
$ws.Range("F268").Value = "Guilty"
$ws.Range("G268").Value = "Guilty"
$ws.Range("H268").Value = "'`$ 0"
$ws.Range("I268").Value = "'`$ 0"

# Row 269
$ws.Range("A269").Value = "21TRD09437"
$ws.Range("B269").Value = "Bunner"
$ws.Range("C269").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D269").Value = "4511.21B1A"
$ws.Range("E269").Value = "M4"
$ws.Range("F269").Value = "Dismissed"
$ws.Range("H269").Value = " "
$ws.Range("I269").Value = " "

# Row 270
$ws.Range("A270").Value = "21TRD09437"
$ws.Range("B270").Value = "Bunner"
$ws.Range("C270").Value = "RECKLESS OPERATION 1ST IN 1 YR"
$ws.Range("D270").Value = "'4511.20"
$ws.Range("E270").Value = "MM"
$ws.Range("F270").Value = "Guilty"
$ws.Range("G270").Value = "Guilty"
$ws.Range("H270").Value = "'`$ 0"
$ws.Range("I270").Value = "'`$ 0"

# Row 271
$ws.Range("A271").Value = "21TRD09437"
$ws.Range("B271").Value = "Bunner"
$ws.Range("C271").Value = "DUS - AMENDED to Distracted Driving"
$ws.Range("D271").Value = "'4510.11"
$ws.Range("E271").Value = "M1"
$ws.Range("F271").Value = "Guilty"
$ws.Range("G271").Value = "Guilty"
$ws.Range("H271").Value = "'`$ 0"
$ws.Range("I271").Value = "'`$ 0"

# Row 272
$ws.Range("A272").Value = "21TRD09437"
$ws.Range("B272").Value = "Bunner"
$ws.Range("C272").Value = "RECKLESS OPERATION 1ST IN 1 YR"
$ws.Range("D272").Value = "'4511.20"
$ws.Range("E272").Value = "MM"
$ws.Range("F272").Value = "Guilty"
$ws.Range("G272").Value = "Guilty"
$ws.Range("H272").Value = "'`$ 0"
$ws.Range("I272").Value = "'`$ 0"

# Row 273
$ws.Range("A273").Value = "21TRD09437"
$ws.Range("B273").Value = "Bunner"
$ws.Range("C273").Value = "DUS"
$ws.Range("D273").Value = "'4510.11"
$ws.Range("E273").Value = "M1"
$ws.Range("F273").Value = "No Contest"
$ws.Range("G273").Value = "Guilty"
$ws.Range("H273").Value = "'`$ 0"
$ws.Range("I273").Value = "'`$ 0"

# Row 274
$ws.Range("A274").Value = "21TRD09437"
$ws.Range("B274").Value = "Bunner"
$ws.Range("C274").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D274").Value = "4511.21B1A"
$ws.Range("E274").Value = "M4"
$ws.Range("F274").Value = "Dismissed"
$ws.Range("H274").Value = " "
$ws.Range("I274").Value = " "

# Row 275
$ws.Range("A275").Value = "21TRD09437"
$ws.Range("B275").Value = "Bunner"
$ws.Range("C275").Value = "RECKLESS OPERATION 1ST IN 1 YR"
$ws.Range("D275").Value = "'4511.20"
$ws.Range("E275").Value = "MM"
$ws.Range("F275").Value = "No Contest"
$ws.Range("G275").Value = "Guilty"
$ws.Range("H275").Value = "'`$ 0"
$ws.Range("I275").Value = "'`$ 0"

# Row 276
$ws.Range("A276").Value = "21TRD09437"
$ws.Range("B276").Value = "Bunner"
$ws.Range("C276").Value = "DUS"
$ws.Range("D276").Value = "'4510.11"
$ws.Range("E276").Value = "M1"
$ws.Range("F276").Value = "Guilty"
$ws.Range("G276").Value = "Guilty"
$ws.Range("H276").Value = "'`$ 0"
$ws.Range("I276").Value = "'`$ 0"
$ws.Range("J276").Value = "None"
$ws.Range("K276").Value = "None"

# Row 277
$ws.Range("A277").Value = "21TRD09437"
$ws.Range("B277").Value = "Bunner"
$ws.Range("C277").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D277").Value = "4511.21B1A"
$ws.Range("E277").Value = "M4"
$ws.Range("F277").Value = "Dismissed"
$ws.Range("H277").Value = " "
$ws.Range("I277").Value = " "
$ws.Range("J277").Value = " "
$ws.Range("K277").Value = " "

# Row 278
$ws.Range("A278").Value = "21TRD09437"
$ws.Range("B278").Value = "Bunner"
$ws.Range("C278").Value = "RECKLESS OPERATION 1ST IN 1 YR"
$ws.Range("D278").Value = "'4511.20"
$ws.Range("E278").Value = "MM"
$ws.Range("F278").Value = "Guilty"
$ws.Range("G278").Value = "Guilty"
$ws.Range("H278").Value = "'`$ 0"
$ws.Range("I278").Value = "'`$ 0"
$ws.Range("J278").Value = "None"
$ws.Range("K278").Value = "None"
